# Added the option to use load curve and distributed generation curve,
# fixed some bugs.
#
# - Adds a new worksheet "LC and DERS" (results for the Load-Curve /
#   Distributed-Energy-Resources-curve scenarios) after the existing
#   "RELRAD" and "MCS" sheets.
# - Makes the new sheet the active sheet/tab.
# - Updates the selection on the RELRAD and MCS sheets to the header row.

$wb = $excel.ActiveWorkbook

$wsRELRAD = $wb.Worksheets.Item(1)
$wsMCS    = $wb.Worksheets.Item(2)

# --- Create the new "LC and DERS" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLC = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsLC.Name = "LC and DERS"

# Header row (B1:H1) - same headers used on the MCS sheet.
$headers = @("SAIFI", "SAIDI", "CAIDI", "EENS", "nr of simulations", "provided beta", "calculated beta")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $wsLC.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Re-use the bold/centered/bordered header formatting already used for the
# header rows on the other two sheets (copy format only, values are untouched).
$wsMCS.Range("B1").Copy() | Out-Null
$wsLC.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Data rows (A2:H6) - one row per scenario. The numeric columns (B:H) are
# filled in row order, but the scenario labels in column A are entered with
# "Base Case" last (matching the order the author originally typed them in)
# so the shared-string table comes out in the same order as the source file.
$rows = @{
    3 = @("DERS",                    1.451406554089917,  9.315717880424746,  6.418406926835209,  48.6627650782175,   2580, 0.02, 0.02079292816962097)
    4 = @("Load Curve",               1.412400051598195,  9.758304734708979,  6.909023207459539,  41.99660502431043,  2287, 0.02, 0.02099958977753318)
    5 = @("Load Curve + DERS",        1.454097759571132,  9.506993309447445,  6.538070254816576,  40.26777210187073,  1872, 0.02, 0.02286345576677238)
    6 = @("Load Curve + DERS Curve",  1.405747720854301,  9.16300767114685,   6.518244728562183,  38.27579855533233,  2241, 0.02, 0.02129529065985877)
    2 = @("Base Case",               1.425855471598278,  9.770998203685073,  6.852726940643225,  51.46480950024928,  2620, 0.02, 0.02088818228051634)
}

foreach ($rowNum in 3, 4, 5, 6, 2) {
    $rowData = $rows[$rowNum]
    $wsLC.Cells.Item($rowNum, 1).Value = $rowData[0]
    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $wsLC.Cells.Item($rowNum, $c + 1).Value = $rowData[$c]
    }
}

# --- Update selections on the existing sheets ---
$wsRELRAD.Range("A1:H1").Select()
$wsMCS.Range("A1:K1").Select()

# --- Make the new sheet the active one, selecting I14 ---
$wsLC.Select()
$wsLC.Range("I14").Select()
